$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row of data
$ws.Range("A3").Value = "Marina"
$ws.Range("B3").Value = "Elswere"

# Move the active selection, matching the saved workbook state
$ws.Range("B7").Select()
